$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "2. Mo ta he thong: db, ngon ngu, "
$ws.Range("D7").Value = "3. Giao dien"
$ws.Range("D5").Value = "1. Hien trang: cac ung dung tuong tu"
